$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44278
$ws.Range("J2").Value2 = 300
$ws.Range("D3").Value2 = 44278
$ws.Range("J3").Value2 = 150
$ws.Range("D4").Value2 = 44194
$ws.Range("O4").Value2 = 'Región de Ñuble'
$ws.Range("D5").Value2 = 44194
$ws.Range("O5").Value2 = 'Región de Ñuble'
$ws.Range("D6").Value2 = 44362
$ws.Range("D7").Value2 = 44362
$ws.Range("D8").Value2 = 44398
$ws.Range("D9").Value2 = 44398
$ws.Range("D10").Value2 = 44222
$ws.Range("D11").Value2 = 44222
$ws.Range("D12").Value2 = 44349
$ws.Range("O12").Value2 = 'Región Metropolitana'
$ws.Range("D13").Value2 = 44349
$ws.Range("O13").Value2 = 'Región Metropolitana'
$ws.Range("D14").Value2 = 44330
$ws.Range("D15").Value2 = 44330
$ws.Range("D16").Value2 = 44231
$ws.Range("D17").Value2 = 44231
$ws.Range("D18").Value2 = 44435
$ws.Range("D19").Value2 = 44435
$ws.Range("D20").Value2 = 44224
$ws.Range("D21").Value2 = 44224
$ws.Range("D22").Value2 = 44313
$ws.Range("D23").Value2 = 44313
$ws.Range("D24").Value2 = 44272
$ws.Range("D25").Value2 = 44272
$ws.Range("D26").Value2 = 44355
$ws.Range("D27").Value2 = 44355
$ws.Range("D28").Value2 = 44299
$ws.Range("D29").Value2 = 44299
$ws.Range("D30").Value2 = 44320
$ws.Range("D31").Value2 = 44320
$ws.Range("D32").Value2 = 44391
$ws.Range("D33").Value2 = 44391
$ws.Range("D34").Value2 = 44334
$ws.Range("D35").Value2 = 44334
$ws.Range("D36").Value2 = 44250
$ws.Range("O36").Value2 = 'Región de Arica y Parinacota'
$ws.Range("D37").Value2 = 44250
$ws.Range("O37").Value2 = 'Región de Arica y Parinacota'
$ws.Range("D38").Value2 = 44327
$ws.Range("D39").Value2 = 44327
$ws.Range("D40").Value2 = 44383
$ws.Range("D41").Value2 = 44383
$ws.Range("D42").Value2 = 44405
$ws.Range("D43").Value2 = 44405
$ws.Range("D44").Value2 = 44336
$ws.Range("N44").Value2 = '$/atado 0,5 a 1 kilo'
$ws.Range("D45").Value2 = 44336
$ws.Range("N45").Value2 = '$/atado 0,5 a 1 kilo'
$ws.Range("D46").Value2 = 44341
$ws.Range("D47").Value2 = 44341
$ws.Range("D48").Value2 = 44274
$ws.Range("J48").Value2 = 200
$ws.Range("D49").Value2 = 44274
$ws.Range("J49").Value2 = 100
$ws.Range("D50").Value2 = 44400
$ws.Range("O50").Value2 = 'Región de Ñuble'
$ws.Range("D51").Value2 = 44400
$ws.Range("O51").Value2 = 'Región de Ñuble'
$ws.Range("D52").Value2 = 44442
$ws.Range("J52").Value2 = 300
$ws.Range("D53").Value2 = 44442
$ws.Range("J53").Value2 = 150
$ws.Range("D54").Value2 = 44453
$ws.Range("D55").Value2 = 44453
$ws.Range("D56").Value2 = 44237
$ws.Range("D57").Value2 = 44237
$ws.Range("D58").Value2 = 44285
$ws.Range("D59").Value2 = 44285
$ws.Range("D60").Value2 = 44344
$ws.Range("N60").Value2 = '$/docena de 1 kilo'
$ws.Range("D61").Value2 = 44344
$ws.Range("N61").Value2 = '$/docena de 1 kilo'
$ws.Range("D62").Value2 = 44217
$ws.Range("D63").Value2 = 44217
$ws.Range("D64").Value2 = 44350
$ws.Range("D65").Value2 = 44350
$ws.Range("D66").Value2 = 44280
$ws.Range("D67").Value2 = 44280
$ws.Range("D68").Value2 = 44447
$ws.Range("D69").Value2 = 44447
$ws.Range("D70").Value2 = 44187
$ws.Range("D71").Value2 = 44187
$ws.Range("D72").Value2 = 44386
$ws.Range("D73").Value2 = 44386
$ws.Range("D74").Value2 = 44308
$ws.Range("D75").Value2 = 44308
$ws.Range("D76").Value2 = 44252
$ws.Range("D77").Value2 = 44252
$ws.Range("D78").Value2 = 44166
$ws.Range("D79").Value2 = 44166
$ws.Range("D80").Value2 = 44168
$ws.Range("D81").Value2 = 44168
$ws.Range("D82").Value2 = 44433
$ws.Range("D83").Value2 = 44433
$ws.Range("D84").Value2 = 44316
$ws.Range("D85").Value2 = 44316
$ws.Range("D86").Value2 = 44160
$ws.Range("D87").Value2 = 44160
$ws.Range("D88").Value2 = 44365
$ws.Range("D89").Value2 = 44365
$ws.Range("D90").Value2 = 44306
$ws.Range("D91").Value2 = 44306
$ws.Range("D92").Value2 = 44203
$ws.Range("D93").Value2 = 44203
$ws.Range("D94").Value2 = 44239
$ws.Range("D95").Value2 = 44239
$ws.Range("D96").Value2 = 44292
$ws.Range("J96").Value2 = 200
$ws.Range("D97").Value2 = 44292
$ws.Range("J97").Value2 = 100
$ws.Range("D98").Value2 = 44358
$ws.Range("D99").Value2 = 44358
$ws.Range("D100").Value2 = 44425
$ws.Range("D101").Value2 = 44425
